$d = $word.ActiveDocument

# Step 1: remove all paragraphs except the last one (cannot delete final paragraph mark)
for ($i = 0; $i -lt 20; $i++) {
  if ($d.Paragraphs.Count -le 1) { break }
  $rng = $d.Content
  $rng.Delete()
}

# Step 2: clear remaining text in the last paragraph, keeping its paragraph mark
$end = $d.Content.End
if ($end -gt 1) {
  $rng = $d.Range(0, $end - 1)
  $rng.Delete()
}

# Step 3: insert the new body content (as WordML) at the start of the document
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Personal Statement</w:t></w:r><w:r><w:t xml:space="preserve"> – David Robertson</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>I am applying</w:t></w:r><w:r><w:t xml:space="preserve"> to study the </w:t></w:r><w:r><w:t>MSc</w:t></w:r><w:r><w:t xml:space="preserve"> in Artificial Intelligence at Aberdeen </w:t></w:r><w:r><w:t>U</w:t></w:r><w:r><w:t>niversity for a number of reasons. Firstly, I feel the content of the course is perfect in terms of what I would like to learn a</w:t></w:r><w:r><w:t>bout particularly machine learning and data mining</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t xml:space="preserve">These areas have interested me greatly since early in my undergraduate studies at </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Abertay</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> University, thus I have based my honours project on the training of a neural network to drive a car in a game. I have also been in contact with the programme coordinator Dr Pang. I had a very interesting conversation with him </w:t></w:r><w:r><w:t>when we met at Data Talent Scotland 2017</w:t></w:r><w:r><w:t xml:space="preserve">, which really increased my interest in the course. </w:t></w:r><w:r><w:t xml:space="preserve">I achieved a very high grade for my research proposal and am on track to continue that trend with my dissertation. I believe I will be able to achieve a first class degree, which will demonstrate my high level of academic ability. </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Through my work with </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Justfone</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> in the summer of 2016, I have increased my ability to work independently with short time frames. I was responsible for producing a python script that would turn a CSV file into a complete web form with links to a database. I had to</w:t></w:r><w:r><w:t xml:space="preserve"> research and</w:t></w:r><w:r><w:t xml:space="preserve"> complete this task within 2 weeks. This helped me develop skills in conducting relevant research and time management.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Through</w:t></w:r><w:r><w:t xml:space="preserve"> the process of executing</w:t></w:r><w:r><w:t xml:space="preserve"> my honours project</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> I have further developed these qualities and increased my understanding of how to conduct</w:t></w:r><w:r><w:t xml:space="preserve"> efficient</w:t></w:r><w:r><w:t xml:space="preserve"> research of certain academic fields.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">I believe this course will expand my knowledge </w:t></w:r><w:r><w:t>in a</w:t></w:r><w:r><w:t xml:space="preserve">rtificial </w:t></w:r><w:r><w:t>i</w:t></w:r><w:r><w:t>ntelligence</w:t></w:r><w:r><w:t xml:space="preserve"> greatly, allowing me to have a much better grasp on the fundamentals of the topic. I</w:t></w:r><w:r><w:t xml:space="preserve"> also</w:t></w:r><w:r><w:t xml:space="preserve"> believ</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve">e </w:t></w:r><w:r><w:t xml:space="preserve">it will </w:t></w:r><w:r><w:t xml:space="preserve">allow me </w:t></w:r><w:r><w:t xml:space="preserve">to pursue further research into the field in the form of a PhD in machine learning. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertPoint = $d.Range(0, 0)
$insertPoint.InsertXML($xml)

# Step 4: remove the now-empty leftover trailing paragraph left over from step 1-2,
# by merging it with the previous paragraph (deleting its paragraph mark).
$n = $d.Paragraphs.Count
if ($n -gt 1) {
  $last = $d.Paragraphs($n)
  $prev = $d.Paragraphs($n - 1)
  $rngDel = $d.Range($prev.Range.End - 1, $last.Range.End)
  $rngDel.Delete()
}

Write-Host "Final paragraph count:" $d.Paragraphs.Count
